$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Free up the three shared strings that are being reshuffled/reworded so the
# new values get appended/reused in the order the final file expects
# (Homepage text first, then the reworded cart/dairy text, then the reworded
# produce text).
$ws.Range("C4").Value = "__tmp_C4__"
$ws.Range("C6").Value = "__tmp_C6__"

$ws.Range("C2").Value = "Homepage P1 (index.html) and everything therein, including the banners, food displays, buttons, and basic styles that are applied elsewhere. Also made P2 grain aisle page and all P3 pages relating to grain foods. I made the user edit page. I also handled hosting of the website and I was the Github administrator, responsible for handling merge conflicts and other issues."

$ws.Range("C4").Value = "Created the shopping cart page (p4), Created the Dairy Aisle and the dairy products (p2 and p3), and created the back page edit product pages (p8)." + [char]10 + " Also worked on CSS relating to the pages created (created a few classes to make the" + [char]10 + " receipt display properly for example)."

$ws.Range("C6").Value = "Created template for P2 pages. Made banners for P2 pages." + [char]10 + " Created Produce aisle (P2), product descriptions for Produce aisle (P3), P5 and P6. Worked on CSS for the pages created."

# New formatting for the reworded cells: centered + wrapped text.
$ws.Range("C4").HorizontalAlignment = -4108
$ws.Range("C4").WrapText = $true

$ws.Range("C6").HorizontalAlignment = -4108
$ws.Range("C6").WrapText = $true

# Explicit row heights to fit the new multi-line text.
$ws.Rows(4).RowHeight = 61.8
$ws.Rows(6).RowHeight = 56.4

# Restore selection to where the author left off.
$ws.Range("C7").Select()
